$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 2.629231666666667
$ws.Cells.Item(2,8).Value = 7.887695
$ws.Cells.Item(2,9).Value = 0.1414315557047068
$ws.Cells.Item(2,10).Value = 0.1414315557047067
$ws.Cells.Item(2,13).Value = 14.69614866666667
$ws.Cells.Item(2,14).Value = 44.088446
$ws.Cells.Item(2,15).Value = 0.7133564251716612
$ws.Cells.Item(2,16).Value = 0.7133564251716613
$ws.Cells.Item(2,17).Value = 38.63957945244111
$ws.Cells.Item(2,18).Value = 347.75621507197
$ws.Cells.Item(2,19).Value = 0.1008911089839763
$ws.Cells.Item(2,20).Value = 0.1008911089839763
$ws.Cells.Item(3,7).Value = 2.629231666666667
$ws.Cells.Item(3,8).Value = 7.887695
$ws.Cells.Item(3,9).Value = 0.1414315557047068
$ws.Cells.Item(3,10).Value = 0.1414315557047067
$ws.Cells.Item(3,15).Value = 0.2169031240449683
$ws.Cells.Item(3,16).Value = 0.2169031240449683
$ws.Cells.Item(3,17).Value = 11.74874887122722
$ws.Cells.Item(3,18).Value = 105.738739841045
$ws.Cells.Item(3,19).Value = 0.03067694627089085
$ws.Cells.Item(3,20).Value = 0.03067694627089085
$ws.Cells.Item(4,7).Value = 2.629231666666667
$ws.Cells.Item(4,8).Value = 7.887695
$ws.Cells.Item(4,9).Value = 0.1414315557047068
$ws.Cells.Item(4,10).Value = 0.1414315557047067
$ws.Cells.Item(4,13).Value = 1.308497333333333
$ws.Cells.Item(4,14).Value = 3.925492
$ws.Cells.Item(4,15).Value = 0.0635149385886714
$ws.Cells.Item(4,16).Value = 0.06351493858867141
$ws.Cells.Item(4,17).Value = 3.440342624548889
$ws.Cells.Item(4,18).Value = 30.96308362094
$ws.Cells.Item(4,19).Value = 0.008983016575084707
$ws.Cells.Item(4,20).Value = 0.008983016575084707
$ws.Cells.Item(5,7).Value = 2.629231666666667
$ws.Cells.Item(5,8).Value = 7.887695
$ws.Cells.Item(5,9).Value = 0.1414315557047068
$ws.Cells.Item(5,10).Value = 0.1414315557047067
$ws.Cells.Item(5,13).Value = 0.1282543333333333
$ws.Cells.Item(5,14).Value = 0.384763
$ws.Cells.Item(5,15).Value = 0.006225512194698899
$ws.Cells.Item(5,16).Value = 0.0062255121946989
$ws.Cells.Item(5,17).Value = 0.3372103545872222
$ws.Cells.Item(5,18).Value = 3.034893191285
$ws.Cells.Item(5,19).Value = 0.0008804838747548886
$ws.Cells.Item(5,20).Value = 0.0008804838747548885
$ws.Cells.Item(6,9).Value = 0.6147160060020365
$ws.Cells.Item(6,10).Value = 0.6147160060020365
$ws.Cells.Item(6,13).Value = 14.69614866666667
$ws.Cells.Item(6,14).Value = 44.088446
$ws.Cells.Item(6,15).Value = 0.7133564251716612
$ws.Cells.Item(6,16).Value = 0.7133564251716613
$ws.Cells.Item(6,17).Value = 167.9424922977956
$ws.Cells.Item(6,18).Value = 1511.48243068016
$ws.Cells.Item(6,19).Value = 0.4385116125374142
$ws.Cells.Item(6,20).Value = 0.4385116125374142
$ws.Cells.Item(7,9).Value = 0.6147160060020365
$ws.Cells.Item(7,10).Value = 0.6147160060020365
$ws.Cells.Item(7,15).Value = 0.2169031240449683
$ws.Cells.Item(7,16).Value = 0.2169031240449683
$ws.Cells.Item(7,19).Value = 0.1333338221022872
$ws.Cells.Item(7,20).Value = 0.1333338221022872
$ws.Cells.Item(8,9).Value = 0.6147160060020365
$ws.Cells.Item(8,10).Value = 0.6147160060020365
$ws.Cells.Item(8,13).Value = 1.308497333333333
$ws.Cells.Item(8,14).Value = 3.925492
$ws.Cells.Item(8,15).Value = 0.0635149385886714
$ws.Cells.Item(8,16).Value = 0.06351493858867141
$ws.Cells.Item(8,17).Value = 14.95305391292445
$ws.Cells.Item(8,18).Value = 134.57748521632
$ws.Cells.Item(8,19).Value = 0.03904364937069271
$ws.Cells.Item(8,20).Value = 0.03904364937069271
$ws.Cells.Item(9,9).Value = 0.6147160060020365
$ws.Cells.Item(9,10).Value = 0.6147160060020365
$ws.Cells.Item(9,13).Value = 0.1282543333333333
$ws.Cells.Item(9,14).Value = 0.384763
$ws.Cells.Item(9,15).Value = 0.006225512194698899
$ws.Cells.Item(9,16).Value = 0.0062255121946989
$ws.Cells.Item(9,17).Value = 1.465646059831111
$ws.Cells.Item(9,18).Value = 13.19081453848
$ws.Cells.Item(9,19).Value = 0.00382692199164228
$ws.Cells.Item(9,20).Value = 0.00382692199164228
$ws.Cells.Item(10,7).Value = 4.24731
$ws.Cells.Item(10,8).Value = 12.74193
$ws.Cells.Item(10,9).Value = 0.2284711798035388
$ws.Cells.Item(10,10).Value = 0.2284711798035388
$ws.Cells.Item(10,13).Value = 14.69614866666667
$ws.Cells.Item(10,14).Value = 44.088446
$ws.Cells.Item(10,15).Value = 0.7133564251716612
$ws.Cells.Item(10,16).Value = 0.7133564251716613
$ws.Cells.Item(10,17).Value = 62.41909919341999
$ws.Cells.Item(10,18).Value = 561.77189274078
$ws.Cells.Item(10,19).Value = 0.1629813840794043
$ws.Cells.Item(10,20).Value = 0.1629813840794043
$ws.Cells.Item(11,7).Value = 4.24731
$ws.Cells.Item(11,8).Value = 12.74193
$ws.Cells.Item(11,9).Value = 0.2284711798035388
$ws.Cells.Item(11,10).Value = 0.2284711798035388
$ws.Cells.Item(11,15).Value = 0.2169031240449683
$ws.Cells.Item(11,16).Value = 0.2169031240449683
$ws.Cells.Item(11,17).Value = 18.97914862387
$ws.Cells.Item(11,18).Value = 170.81233761483
$ws.Cells.Item(11,19).Value = 0.04955611265362723
$ws.Cells.Item(11,20).Value = 0.04955611265362724
$ws.Cells.Item(12,7).Value = 4.24731
$ws.Cells.Item(12,8).Value = 12.74193
$ws.Cells.Item(12,9).Value = 0.2284711798035388
$ws.Cells.Item(12,10).Value = 0.2284711798035388
$ws.Cells.Item(12,13).Value = 1.308497333333333
$ws.Cells.Item(12,14).Value = 3.925492
$ws.Cells.Item(12,15).Value = 0.0635149385886714
$ws.Cells.Item(12,16).Value = 0.06351493858867141
$ws.Cells.Item(12,17).Value = 5.55759380884
$ws.Cells.Item(12,18).Value = 50.01834427956
$ws.Cells.Item(12,19).Value = 0.01451133295450307
$ws.Cells.Item(12,20).Value = 0.01451133295450307
$ws.Cells.Item(13,7).Value = 4.24731
$ws.Cells.Item(13,8).Value = 12.74193
$ws.Cells.Item(13,9).Value = 0.2284711798035388
$ws.Cells.Item(13,10).Value = 0.2284711798035388
$ws.Cells.Item(13,13).Value = 0.1282543333333333
$ws.Cells.Item(13,14).Value = 0.384763
$ws.Cells.Item(13,15).Value = 0.006225512194698899
$ws.Cells.Item(13,16).Value = 0.0062255121946989
$ws.Cells.Item(13,17).Value = 0.54473591251
$ws.Cells.Item(13,18).Value = 4.90262321259
$ws.Cells.Item(13,19).Value = 0.001422350116004176
$ws.Cells.Item(13,20).Value = 0.001422350116004176
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.2859396666666667
$ws.Cells.Item(14,8).Value = 0.8578190000000001
$ws.Cells.Item(14,9).Value = 0.01538125848971795
$ws.Cells.Item(14,10).Value = 0.01538125848971795
$ws.Cells.Item(14,13).Value = 14.69614866666667
$ws.Cells.Item(14,14).Value = 44.088446
$ws.Cells.Item(14,15).Value = 0.7133564251716612
$ws.Cells.Item(14,16).Value = 0.7133564251716613
$ws.Cells.Item(14,17).Value = 4.202211851030444
$ws.Cells.Item(14,18).Value = 37.819906659274
$ws.Cells.Item(14,19).Value = 0.01097231957086646
$ws.Cells.Item(14,20).Value = 0.01097231957086646
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.2859396666666667
$ws.Cells.Item(15,8).Value = 0.8578190000000001
$ws.Cells.Item(15,9).Value = 0.01538125848971795
$ws.Cells.Item(15,10).Value = 0.01538125848971795
$ws.Cells.Item(15,15).Value = 0.2169031240449683
$ws.Cells.Item(15,16).Value = 0.2169031240449683
$ws.Cells.Item(15,17).Value = 1.277724355209889
$ws.Cells.Item(15,18).Value = 11.499519196889
$ws.Cells.Item(15,19).Value = 0.003336243018163015
$ws.Cells.Item(15,20).Value = 0.003336243018163015
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.2859396666666667
$ws.Cells.Item(16,8).Value = 0.8578190000000001
$ws.Cells.Item(16,9).Value = 0.01538125848971795
$ws.Cells.Item(16,10).Value = 0.01538125848971795
$ws.Cells.Item(16,13).Value = 1.308497333333333
$ws.Cells.Item(16,14).Value = 3.925492
$ws.Cells.Item(16,15).Value = 0.0635149385886714
$ws.Cells.Item(16,16).Value = 0.06351493858867141
$ws.Cells.Item(16,17).Value = 0.3741512913275556
$ws.Cells.Item(16,18).Value = 3.367361621948
$ws.Cells.Item(16,19).Value = 0.0009769396883909163
$ws.Cells.Item(16,20).Value = 0.0009769396883909166
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.2859396666666667
$ws.Cells.Item(17,8).Value = 0.8578190000000001
$ws.Cells.Item(17,9).Value = 0.01538125848971795
$ws.Cells.Item(17,10).Value = 0.01538125848971795
$ws.Cells.Item(17,13).Value = 0.1282543333333333
$ws.Cells.Item(17,14).Value = 0.384763
$ws.Cells.Item(17,15).Value = 0.006225512194698899
$ws.Cells.Item(17,16).Value = 0.0062255121946989
$ws.Cells.Item(17,17).Value = 0.03667300132188889
$ws.Cells.Item(17,18).Value = 0.330057011897
$ws.Cells.Item(17,19).Value = 0.00009575621229755509
$ws.Cells.Item(17,20).Value = 0.0000957562122975551
